$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update metric values (B2:B13) ---
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 61347.969999999994
$metrics.Range("B3").Value = 49873.91
$metrics.Range("B4").Value = 19421.170000000002
$metrics.Range("B5").Value = 2435
$metrics.Range("B6").Value = 3980598.8499999996
$metrics.Range("B7").Value = 3377401.3899999992
$metrics.Range("B8").Value = 1148786.8499999999
$metrics.Range("B9").Value = 153595
$metrics.Range("B10").Value = 32445922.650999829
$metrics.Range("B11").Value = 19407271.460000005
$metrics.Range("B12").Value = 11430495.740000002
$metrics.Range("B13").Value = 1251222

# Move the Metrics selection to D25 (was D18)
$metrics.Range("D25").Select()

# --- zgmysj sheet: move selection from B78:B85 to P90 (no longer the active tab) ---
$zgmysj = $wb.Worksheets.Item("zgmysj")
$zgmysj.Range("P90").Select()

# --- today sheet: move selection from H11 to D5; ends up the active tab ---
$today = $wb.Worksheets.Item("today")
$today.Range("D5").Select()
